# Error Calculations and Plots
# Apply cell-level corrections (values filled in / cleared) plus removal of
# two outlier rows (RM 232 and SC 92), which shifts all subsequent rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: cell-level edits that use the ORIGINAL row numbering
# (rows 2-25 sit above the two rows that get deleted below, so their row
# numbers are unaffected by that later shift).
$ws.Range("C3").Value = 11.2
$ws.Range("D4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("E6").Value = -5.7
$ws.Range("D9").Value = -14.5
$ws.Range("D10").Value = -14.7
$ws.Range("E12").ClearContents()
$ws.Range("E14").Value = -5.4
$ws.Range("D17").ClearContents()
$ws.Range("E17").Value = -7.3
$ws.Range("D18").ClearContents()
$ws.Range("E19").Value = -6.5
$ws.Range("E20").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("E25").ClearContents()

# --- Step 2: remove the two outlier data rows entirely (RM 232 and SC 92).
# Delete from the bottom up so the earlier deletion's row index stays valid.
$ws.Rows.Item(28).Delete()   # "SC 92" row
$ws.Rows.Item(26).Delete()   # "RM 232" row

# --- Step 3: cell-level edits that only make sense using the NEW row
# numbering, i.e. for rows that shifted up after the deletions above
# (originally SC 101, SC 105, SC 193 -> now rows 27, 28, 32).
$ws.Range("E27").Value = -10
$ws.Range("E28").Value = -5.9
$ws.Range("C32").Value = 10.5
